$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values in column D and percentage strings in column E
# stay stored as text (matching the original inlineStr cell type) instead of being
# auto-converted to numbers by Excel when assigned through .Value.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "91.964.76"
$ws.Range("E2").Value = "  +5.02%  "
$ws.Range("D3").Value = "3.268.88"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "217.88"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").Value = "628.43"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "0.413"
$ws.Range("E7").Value = "  +9.27%  "
$ws.Range("D8").Value = "0.719"
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "3.267.96"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Value = "0.588"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").Value = "34.30"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "3.869.48"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "91.804.93"
$ws.Range("E16").Value = "  +5.11%  "
$ws.Range("D17").Value = "5.35"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "3.263.84"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +7.10%  "
$ws.Range("D20").Value = "14.02"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "439.52"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("E22").Value = "  +49.95%  "
$ws.Range("D23").Value = "8.89"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "5.27"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("E25").Value = "  +5.22%  "
$ws.Range("D26").Value = "12.38"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "3.450.73"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "77.15"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "0.179"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "8.78"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "554.50"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "7.16"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").Value = "3.65"
$ws.Range("E35").Value = "  +24.48%  "
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "1.93"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").Value = "1.29"
$ws.Range("E37").Value = "  -7.43%  "
$ws.Range("D38").Value = "22.65"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  -4.75%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "150.09"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("D46").Value = "179.32"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").Value = "45.14"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("E48").Value = "  +5.87%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "4.24"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E51").Value = "  +2.32%  "

# Restore the default (unstyled) cell style now that the values are stored as text,
# so the resulting style indices match the original workbook formatting.
$dataRange.Style = "Normal"

